# Alaska_Review.xlsx - "Updated review sheet 0.3.1.16b"
# Adds a new "Appended" column (H) with Yes/No values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Header (H1): same style as the other header cells (A1:G1) ---
$ws.Range("H1").Value = "Appended"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial($xlPasteFormats)

# --- Data rows: Yes (green fill) / No (red fill, white font like col F "Incomplete") ---
$yesRows = @(2,3,4,6,8,9,10,11,12,18)
$noRows  = @(5,7,13,14,15,16,17,19,20,21,22,23)

foreach ($r in $yesRows) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Value = "Yes"
    $cell.Interior.Color = 5296274
}

foreach ($r in $noRows) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Value = "No"
}
# Paint the "No" cells with the same format already used for "Incomplete" (F11)
$ws.Range("F11").Copy()
foreach ($r in $noRows) {
    $ws.Cells.Item($r, 8).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

# --- Column widths: G and H both end up the same width as the old H ---
$ws.Columns("G:H").ColumnWidth = 12.33

# --- Selection moved to G6 ---
$ws.Range("G6").Select()
